$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current header row (row 2) is: B=User, C=Email, D=Topics, E=Regions, F=Alert Slots
# We need to insert two new columns ("Excluded Topics", "Real time Topics") right after
# "Topics" (column D), pushing "Regions" and "Alert Slots" two columns to the right
# (columns G and H).

# Move the existing "Regions"/"Alert Slots" header values two columns to the right.
$ws.Range("H2").Value = $ws.Range("F2").Text
$ws.Range("G2").Value = $ws.Range("E2").Text

# Fill the freed-up columns with the two new headers.
$ws.Range("E2").Value = "Excluded Topics"
$ws.Range("F2").Value = "Real time Topics"

# Match the bold header formatting used by the rest of the row.
$ws.Range("G2:H2").Font.Bold = $true

# Update the selection to cover the new header extent.
$ws.Range("B2:H2").Select()
